$d = $word.ActiveDocument

$replacements = @(
    @("23÷5=", "67÷5="),
    @("81÷6=", "88÷9="),
    @("28÷2=", "56÷5="),
    @("69÷6=", "41÷8="),
    @("74÷9=", "63÷9="),
    @("64÷4=", "66÷3="),
    @("59÷6=", "12÷9="),
    @("46÷7=", "28÷9="),
    @("11÷9=", "37÷6="),
    @("86÷3=", "31÷5="),
    @("45÷4=", "40÷7="),
    @("46÷4=", "47÷6="),
    @("37÷8=", "92÷7="),
    @("42÷2=", "79÷7="),
    @("98÷5=", "54÷7="),
    @("85÷3=", "65÷3="),
    @("84÷9=", "51÷4="),
    @("68÷5=", "78÷5="),
    @("29÷4=", "40÷2="),
    @("22÷3=", "73÷3="),
    @("13÷7=", "63÷7="),
    @("50÷6=", "59÷7="),
    @("81÷7=", "97÷4="),
    @("41÷6=", "88÷4="),
    @("53÷4=", "91÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
